{"js": "// The document's single table holds 20 rows x 5 cols of \"a op b = c\"\n// equations, one equation per table cell. The commit replaces each of\n// the 100 equations with a new one while leaving everything else in the\n// OOXML (paragraph justification, run fonts/size, table structure, etc.)\n// untouched. We walk the table cells in row-major order and overwrite\n// just the text of each cell's range, which keeps the existing\n// run/paragraph formatting intact (unlike replacing the whole cell body,\n// which would drop the paragraph/run properties).\n\nconst newValues = [\n  \"51-23=28\", \"12+66=78\", \"21+56=77\", \"42-16=26\", \"85-23=62\",\n  \"15+71=86\", \"53+0=53\", \"20+2=22\", \"96-1=95\", \"26+16=42\",\n  \"28+16=44\", \"92-92=0\", \"43-39=4\", \"77+15=92\", \"39-38=1\",\n  \"10+41=51\", \"10+52=62\", \"91-13=78\", \"64+17=81\", \"59+21=80\",\n  \"31+25=56\", \"48-15=33\", \"8+17=25\", \"56-11=45\", \"41+7=48\",\n  \"68+20=88\", \"37+61=98\", \"77+14=91\", \"83-51=32\", \"72-12=60\",\n  \"2+63=65\", \"99-0=99\", \"7+62=69\", \"83+6=89\", \"45+46=91\",\n  \"0+97=97\", \"21-2=19\", \"68-20=48\", \"52-50=2\", \"65-23=42\",\n  \"82-19=63\", \"16+42=58\", \"61-21=40\", \"18+58=76\", \"3+16=19\",\n  \"38+29=67\", \"28-13=15\", \"68-17=51\", \"77-22=55\", \"24+30=54\",\n  \"67-1=66\", \"25+57=82\", \"69+2=71\", \"49-24=25\", \"59+4=63\",\n  \"76+0=76\", \"31+58=89\", \"55-13=42\", \"49+32=81\", \"73-55=18\",\n  \"13+39=52\", \"96-51=45\", \"39-19=20\", \"70-30=40\", \"81-61=20\",\n  \"0+60=60\", \"60-25=35\", \"58-0=58\", \"87-85=2\", \"77-9=68\",\n  \"53-23=30\", \"97-15=82\", \"43+33=76\", \"84+8=92\", \"30+42=72\",\n  \"26+31=57\", \"84-45=39\", \"42-38=4\", \"92-75=17\", \"54-20=34\",\n  \"57-6=51\", \"25+11=36\", \"98-27=71\", \"70-15=55\", \"15-6=9\",\n  \"46-45=1\", \"30+20=50\", \"62-6=56\", \"48-13=35\", \"26+70=96\",\n  \"56-14=42\", \"43-11=32\", \"16+41=57\", \"72-51=21\", \"65-3=62\",\n  \"3+53=56\", \"50-37=13\", \"15+62=77\", \"45-41=4\", \"2+68=70\"\n];\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst rows = table.rowCount;\nconst cols = table.values[0].length;\n\nlet i = 0;\nfor (let r = 0; r < rows; r++) {\n  for (let c = 0; c < cols; c++) {\n    const cell = table.getCell(r, c);\n    const range = cell.body.getRange();\n    range.insertText(newValues[i], Word.InsertLocation.replace);\n    i++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document's single table holds 20 rows x 5 cols of \"a op b = c\"\n# equations, one equation per table cell. The commit replaces each of\n# the 100 equations with a new one while leaving every other part of\n# the OOXML (paragraph justification, run fonts/size, table structure,\n# etc.) untouched. We walk the table cells in row-major order (matching\n# the order the equations appear in the source document) and overwrite\n# just each cell's Range.Text, which keeps the existing run/paragraph\n# formatting intact.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValues = @(\n    \"51-23=28\",\"12+66=78\",\"21+56=77\",\"42-16=26\",\"85-23=62\",\n    \"15+71=86\",\"53+0=53\",\"20+2=22\",\"96-1=95\",\"26+16=42\",\n    \"28+16=44\",\"92-92=0\",\"43-39=4\",\"77+15=92\",\"39-38=1\",\n    \"10+41=51\",\"10+52=62\",\"91-13=78\",\"64+17=81\",\"59+21=80\",\n    \"31+25=56\",\"48-15=33\",\"8+17=25\",\"56-11=45\",\"41+7=48\",\n    \"68+20=88\",\"37+61=98\",\"77+14=91\",\"83-51=32\",\"72-12=60\",\n    \"2+63=65\",\"99-0=99\",\"7+62=69\",\"83+6=89\",\"45+46=91\",\n    \"0+97=97\",\"21-2=19\",\"68-20=48\",\"52-50=2\",\"65-23=42\",\n    \"82-19=63\",\"16+42=58\",\"61-21=40\",\"18+58=76\",\"3+16=19\",\n    \"38+29=67\",\"28-13=15\",\"68-17=51\",\"77-22=55\",\"24+30=54\",\n    \"67-1=66\",\"25+57=82\",\"69+2=71\",\"49-24=25\",\"59+4=63\",\n    \"76+0=76\",\"31+58=89\",\"55-13=42\",\"49+32=81\",\"73-55=18\",\n    \"13+39=52\",\"96-51=45\",\"39-19=20\",\"70-30=40\",\"81-61=20\",\n    \"0+60=60\",\"60-25=35\",\"58-0=58\",\"87-85=2\",\"77-9=68\",\n    \"53-23=30\",\"97-15=82\",\"43+33=76\",\"84+8=92\",\"30+42=72\",\n    \"26+31=57\",\"84-45=39\",\"42-38=4\",\"92-75=17\",\"54-20=34\",\n    \"57-6=51\",\"25+11=36\",\"98-27=71\",\"70-15=55\",\"15-6=9\",\n    \"46-45=1\",\"30+20=50\",\"62-6=56\",\"48-13=35\",\"26+70=96\",\n    \"56-14=42\",\"43-11=32\",\"16+41=57\",\"72-51=21\",\"65-3=62\",\n    \"3+53=56\",\"50-37=13\",\"15+62=77\",\"45-41=4\",\"2+68=70\"\n)\n\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\n$i = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $newValues[$i]\n        $i = $i + 1\n    }\n}\n"}
